# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K" = strikeouts) on the active sheet with the
# freshly-regenerated per-game strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 3
    9  = 0
    10 = 2
    12 = 0
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 2
    20 = 2
    21 = 1
    22 = 0
    23 = 1
    24 = 0
    25 = 1
    26 = 0
    27 = 2
    28 = 1
    29 = 1
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 1
    35 = 0
    36 = 2
    37 = 0
    38 = 1
    39 = 3
    40 = 1
    41 = 1
    42 = 0
    44 = 1
    46 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
